# Apply the edit described by the diff:
#  - On "safety_orders": delete row 2 (the old Safety Order No. 3 row),
#    shifting all subsequent rows up (old row 3 -> new row 2, etc.)
#  - On "open_buy_orders": append a new row (row 4) with a new txid/price pair

$wb = $excel.ActiveWorkbook

# --- Sheet: safety_orders ---
$wsSafety = $wb.Worksheets.Item("safety_orders")
$wsSafety.Rows.Item(2).Delete()

# --- Sheet: open_buy_orders ---
$wsBuy = $wb.Worksheets.Item("open_buy_orders")
$wsBuy.Cells.Item(4, 1).Value = "OIBWEE-4TLOO-BFC25Z"
$wsBuy.Cells.Item(4, 2).Value = 6.1311
